$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the PaymentType label "Existing Liability w/Notice Number" to the
# new, more descriptive "Existing Liability with Notice/Invoice Number" text
# everywhere it's used in column D.
$ws.Range("D1:D30").Replace("Existing Liability w/Notice Number", "Existing Liability with Notice/Invoice Number")

# Fill in the (previously mostly-empty) Execute column C with "Y" for every
# data row, matching rows 11-12 that already had it.
for ($r = 2; $r -le 30; $r++) {
    $ws.Cells.Item($r, 3).Value = "Y"
}

# Update the active selection to reflect the newly filled column C range.
$ws.Range("C2:C30").Select()
